# Update the "EC" (Estado de Cuenta) worksheet:
#  - Replace the two mora periods (2507, 2506) with a single new period (2508)
#  - Remove the now-obsolete rows for the second period, keeping one row per worker
#  - Update the totals (VALOR MORA) and period count accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows that represented the second (now-removed) period (rows 17 and 18),
# which shifts row 19 up to row 17 and the footer rows up accordingly.
$ws.Rows.Item(17).Resize(2).Delete()

# Update the period label for both remaining worker rows (now rows 16 and 17) to "2508"
$ws.Range("E16").Value = "2508"
$ws.Range("E17").Value = "2508"

# Update the aggregate "VALOR MORA" total (was 227760 for 2 periods, now for 1 period)
$ws.Range("E11").Value = 113880

# Update "Cant. Periodos" count (was 2, now 1)
$ws.Range("F13").Value = 1
